$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/CodeSystem/fr-medication-reconciliation-type"
$ws.Range("B4").Value = "FRMedicationReconciliationType"
$ws.Range("B5").Value = "code system Interop'Santé - Type d'écart/erreur sur une ligne de traitement d'une FCT"
$ws.Range("B8").Value = "2026-01-15T08:54:26+00:00"
$ws.Range("B11").Value = "FRANCE"
